# Update Name of Algo
# Apply updated numeric values to result_data_RandomForest worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.5831
$ws.Range("A3").Value = -21.40810000000003
$ws.Range("C5").Value = -14.4477
$ws.Range("D5").Value = -8.543799999999997
$ws.Range("D9").Value = -8.051600000000001
$ws.Range("D11").Value = -8.344200000000003
$ws.Range("A14").Value = -20.4423
$ws.Range("A16").Value = -20.2064
$ws.Range("C16").Value = -11.63429999999999
$ws.Range("D17").Value = -8.314400000000004
$ws.Range("A21").Value = -21.13430000000001
$ws.Range("D21").Value = -7.708500000000003
$ws.Range("A23").Value = -21.70900000000003
$ws.Range("A25").Value = -22.45680000000003
